$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 40
$ws.Range("I12").Value = 40
$ws.Range("K12").Value = 40
$ws.Range("M12").Value = 130

$ws.Range("H15").Value = 1066
$ws.Range("I15").Value = 1066
$ws.Range("K15").Value = 3198
$ws.Range("M15").Value = -3029

$ws.Range("H18").Value = 1039.8
$ws.Range("J18").Value = 1408.3334
$ws.Range("L18").Value = 1408.3334
$ws.Range("N18").Value = -1976.3334

$ws.Range("H20").Value = 12470.25
$ws.Range("I20").Value = 12470.25
$ws.Range("K20").Value = 12470.25
$ws.Range("M20").Value = -12240.25

$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1674

$ws.Range("H33").Value = 269.45
$ws.Range("I33").Value = 223.23529
$ws.Range("J33").Value = 531.3333
$ws.Range("K33").Value = 223.23529
$ws.Range("L33").Value = 531.3333
$ws.Range("M33").Value = 5.764710000000008
$ws.Range("N33").Value = -989.3333

$ws.Range("H35").Value = 12470.25
$ws.Range("I35").Value = 12470.25
$ws.Range("K35").Value = 12470.25
$ws.Range("M35").Value = -12091.25

$ws.Range("H108").Value = 99899.42999999999
$ws.Range("J108").Value = 99899.42999999999
$ws.Range("L108").Value = 99899.42999999999
$ws.Range("N108").Value = -107579.43

$ws.Range("H116").Value = 6920.222
$ws.Range("I116").Value = 6345.4546
$ws.Range("K116").Value = 6345.4546
$ws.Range("M116").Value = -2903.4546

$ws.Range("H138").Value = 1492.2142
$ws.Range("J138").Value = 1917.45
$ws.Range("L138").Value = 5752.35
$ws.Range("N138").Value = -16032.35

$ws.Range("H141").Value = 8336.5
$ws.Range("I141").Value = 7505
$ws.Range("K141").Value = 22515
$ws.Range("M141").Value = -17335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7300.0586
$ws.Range("I32").Value = 3264.6875
$ws.Range("J32").Value = 14096.474
$ws.Range("K32").Value = 3264.6875
$ws.Range("L32").Value = 14096.474
$ws.Range("M32").Value = -2977.6875
$ws.Range("N32").Value = -14670.474

$ws.Range("H74").Value = 2049.647
$ws.Range("I74").Value = 1464.5454
$ws.Range("K74").Value = 1464.5454
$ws.Range("M74").Value = -590.5454

$ws.Range("H77").Value = 2049.647
$ws.Range("I77").Value = 1464.5454
$ws.Range("K77").Value = 7322.727
$ws.Range("M77").Value = -2954.727

$ws.Range("H122").Value = 1917.2084
$ws.Range("I122").Value = 1899.9474
$ws.Range("K122").Value = 5699.8422
$ws.Range("M122").Value = -3249.8422

$ws.Range("H134").Value = 124428.5
$ws.Range("J134").Value = 124428.5
$ws.Range("L134").Value = 124428.5
$ws.Range("N134").Value = -134568.5

$ws.Range("H135").Value = 32166.334
$ws.Range("J135").Value = 32166.334
$ws.Range("L135").Value = 32166.334
$ws.Range("N135").Value = -42306.334

$ws.Range("H138").Value = 62763.8
$ws.Range("J138").Value = 59857.25
$ws.Range("L138").Value = 59857.25
$ws.Range("N138").Value = -70137.25

$ws.Range("H140").Value = 70995.75
$ws.Range("J140").Value = 70995.75
$ws.Range("L140").Value = 70995.75
$ws.Range("N140").Value = -81355.75

$ws.Range("H141").Value = 92425.42999999999
$ws.Range("J141").Value = 84195.60000000001
$ws.Range("L141").Value = 84195.60000000001
$ws.Range("N141").Value = -94555.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 60662.766
$ws.Range("I99").Value = 112164.11
$ws.Range("K99").Value = 112164.11
$ws.Range("M99").Value = -110666.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1718.2778
$ws.Range("I58").Value = 1528.8667
$ws.Range("K58").Value = 1528.8667
$ws.Range("M58").Value = -1325.8667

$ws.Range("H134").Value = 3854.3125
$ws.Range("I134").Value = 3940.6428
$ws.Range("K134").Value = 11821.9284
$ws.Range("M134").Value = -9286.928400000001

$ws.Range("H136").Value = 1718.2778
$ws.Range("I136").Value = 1528.8667
$ws.Range("K136").Value = 4586.6001
$ws.Range("M136").Value = -2036.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.733334
$ws.Range("I2").Value = 22.666666
$ws.Range("J2").Value = 50.333332
$ws.Range("K2").Value = 135.999996
$ws.Range("L2").Value = 301.999992
$ws.Range("M2").Value = -22.99999600000001
$ws.Range("N2").Value = -527.999992

$ws.Range("H107").Value = 1747.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1747.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5242.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9082.5

$ws.Range("H128").Value = 381321
$ws.Range("I128").Value = 381321
$ws.Range("K128").Value = 1143963
$ws.Range("M128").Value = -1138983

$ws.Range("H132").Value = 5964.95
$ws.Range("J132").Value = 7965.2144
$ws.Range("L132").Value = 71686.9296
$ws.Range("N132").Value = -76746.9296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14509.143
$ws.Range("J26").Value = 14509.143
$ws.Range("L26").Value = 14509.143
$ws.Range("N26").Value = -15069.143

$ws.Range("H43").Value = 25003154
$ws.Range("I43").Value = 25003154
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 25003154
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -25003003
$ws.Range("N43").ClearContents()

$ws.Range("H50").Value = 14509.143
$ws.Range("J50").Value = 14509.143
$ws.Range("L50").Value = 14509.143
$ws.Range("N50").Value = -15505.143

$ws.Range("H97").Value = 1395.6666
$ws.Range("I97").Value = 668.6111
$ws.Range("J97").Value = 5758
$ws.Range("K97").Value = 668.6111
$ws.Range("L97").Value = 5758
$ws.Range("M97").Value = -172.6111
$ws.Range("N97").Value = -6750

$ws.Range("H102").Value = 1615
$ws.Range("I102").Value = 1615
$ws.Range("K102").Value = 1615
$ws.Range("M102").Value = 7

$ws.Range("H138").Value = 150000
$ws.Range("J138").Value = 150000
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

$ws.Range("H141").Value = 83104.125
$ws.Range("J141").Value = 88805.5
$ws.Range("L141").Value = 88805.5
$ws.Range("N141").Value = -99165.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3573.158
$ws.Range("I46").Value = 2375
$ws.Range("K46").Value = 2375
$ws.Range("M46").Value = -2187

$ws.Range("H55").Value = 2186.6
$ws.Range("I55").Value = 868.3333
$ws.Range("J55").Value = 4164
$ws.Range("K55").Value = 868.3333
$ws.Range("L55").Value = 4164
$ws.Range("M55").Value = -695.3333
$ws.Range("N55").Value = -4510

$ws.Range("H134").Value = 94797
$ws.Range("J134").Value = 94797
$ws.Range("L134").Value = 94797
$ws.Range("N134").Value = -104937

$ws.Range("H135").Value = 76329.664
$ws.Range("J135").Value = 76329.664
$ws.Range("L135").Value = 76329.664
$ws.Range("N135").Value = -86469.664

$ws.Range("H138").Value = 92444.664
$ws.Range("J138").Value = 92444.664
$ws.Range("L138").Value = 92444.664
$ws.Range("N138").Value = -102724.664

$ws.Range("H140").Value = 69799.5
$ws.Range("J140").Value = 69799.5
$ws.Range("L140").Value = 69799.5
$ws.Range("N140").Value = -80159.5

$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 70399.89
$ws.Range("J46").Value = 70399.89
$ws.Range("L46").Value = 70399.89
$ws.Range("N46").Value = -70861.89

$ws.Range("H134").Value = 70399.89
$ws.Range("J134").Value = 70399.89
$ws.Range("L134").Value = 211199.67
$ws.Range("N134").Value = -216269.67

$ws.Range("H135").Value = 63082.332
$ws.Range("J135").Value = 63082.332
$ws.Range("L135").Value = 63082.332
$ws.Range("N135").Value = -73222.33199999999

$ws.Range("H137").Value = 134571.28
$ws.Range("J137").Value = 134571.28
$ws.Range("L137").Value = 134571.28
$ws.Range("N137").Value = -144771.28

$ws.Range("H138").Value = 81693.336
$ws.Range("J138").Value = 81693.336
$ws.Range("L138").Value = 81693.336
$ws.Range("N138").Value = -91973.336

$ws.Range("H140").Value = 124345.43
$ws.Range("J140").Value = 124345.43
$ws.Range("L140").Value = 124345.43
$ws.Range("N140").Value = -134705.43

$ws.Range("H141").Value = 60902.2
$ws.Range("J141").Value = 60902.2
$ws.Range("L141").Value = 60902.2
$ws.Range("N141").Value = -71262.2
